$ws = $excel.ActiveWorkbook.ActiveSheet

# --- Header text updates (Volume/Number and week-covering dates) ---
$ws.Range("A8").Value = "Volume 32   Number  42"
$ws.Range("C9").Value = "Report Covering the Week  10/13/2025  Through  10/19/2025"

# --- Pure numeric value updates (format/type unchanged) ---
$ws.Range("J15").Value = 6
$ws.Range("K15").Value = 0
$ws.Range("J16").Value = 20
$ws.Range("K16").Value = 0
$ws.Range("N16").Value = -59.183673469387
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 6
$ws.Range("H17").Value = -40
$ws.Range("I17").Value = 85
$ws.Range("J17").Value = 77
$ws.Range("K17").Value = 10.389610389610
$ws.Range("L17").Value = 16.438356164383
$ws.Range("M17").Value = 129.72972972973
$ws.Range("N17").Value = -6.593406593406
$ws.Range("L18").Value = -17.391304347826
$ws.Range("M18").Value = -56.818181818181
$ws.Range("N18").Value = -85.767790262172
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 100
$ws.Range("G19").Value = 14
$ws.Range("H19").Value = 42.857142857142
$ws.Range("I19").Value = 238
$ws.Range("J19").Value = 225
$ws.Range("K19").Value = 5.777777777777
$ws.Range("L19").Value = 3.478260869565
$ws.Range("M19").Value = 108.771929824561
$ws.Range("N19").Value = 27.272727272727
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 21
$ws.Range("K20").Value = -38.235294117647
$ws.Range("L20").Value = -67.692307692307
$ws.Range("M20").Value = -32.258064516129
$ws.Range("N20").Value = -96.360485268630
$ws.Range("C21").Value = 8
$ws.Range("D21").Value = 5
$ws.Range("E21").Value = 60
$ws.Range("G21").Value = 34
$ws.Range("H21").Value = -5.882352941176
$ws.Range("I21").Value = 408
$ws.Range("J21").Value = 396
$ws.Range("K21").Value = 3.030303030303
$ws.Range("L21").Value = -5.555555555555
$ws.Range("M21").Value = 37.837837837837
$ws.Range("N21").Value = -65.276595744680
$ws.Range("C24").Value = 5
$ws.Range("D24").Value = 4
$ws.Range("E24").Value = 25
$ws.Range("F24").Value = 40
$ws.Range("G24").Value = 26
$ws.Range("H24").Value = 53.846153846153
$ws.Range("I24").Value = 357
$ws.Range("J24").Value = 291
$ws.Range("K24").Value = 22.680412371134
$ws.Range("L24").Value = -7.512953367875
$ws.Range("M24").Value = -19.413092550790
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 17
$ws.Range("G25").Value = 10
$ws.Range("H25").Value = 70
$ws.Range("I25").Value = 203
$ws.Range("J25").Value = 129
$ws.Range("K25").Value = 57.364341085271
$ws.Range("L25").Value = 15.340909090909
$ws.Range("C26").Value = 7
$ws.Range("E26").Value = 133.333333333333
$ws.Range("F26").Value = 21
$ws.Range("G26").Value = 14
$ws.Range("H26").Value = 50
$ws.Range("I26").Value = 173
$ws.Range("J26").Value = 132
$ws.Range("K26").Value = 31.060606060606
$ws.Range("L26").Value = 7.453416149068
$ws.Range("M26").Value = -3.888888888888
$ws.Range("J27").Value = 10
$ws.Range("K27").Value = -30

# --- Cells switching from numeric to text placeholder ("0" / "***.*") ---
# NumberFormat must be forced to text ("@") BEFORE assigning the value,
# otherwise a numeric-looking string like "0" is coerced back to a number.
# Afterwards the format is switched to "General" to match the target cell style.
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "0"
$ws.Range("C16").NumberFormat = "General"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "0"
$ws.Range("C18").NumberFormat = "General"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "***.*"
$ws.Range("E18").NumberFormat = "General"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "***.*"
$ws.Range("E20").NumberFormat = "General"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "***.*"
$ws.Range("E28").NumberFormat = "General"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "0"
$ws.Range("G31").NumberFormat = "General"
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "***.*"
$ws.Range("H31").NumberFormat = "General"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "***.*"
$ws.Range("E33").NumberFormat = "General"

# --- Cells switching from text placeholder to numeric ---
$ws.Range("M14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M14").Value = -100
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("D15").Value = 1
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E15").Value = -100
$ws.Range("G15").NumberFormat = '#,##0'
$ws.Range("G15").Value = 1
$ws.Range("H15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H15").Value = -100
$ws.Range("D16").NumberFormat = '#,##0'
$ws.Range("D16").Value = 1
$ws.Range("E16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E16").Value = -100
$ws.Range("G16").NumberFormat = '#,##0'
$ws.Range("G16").Value = 1
$ws.Range("H16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H16").Value = 0
$ws.Range("C17").NumberFormat = '#,##0'
$ws.Range("C17").Value = 2
$ws.Range("C20").NumberFormat = '#,##0'
$ws.Range("C20").Value = 2
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("D27").Value = 1
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E27").Value = -100
$ws.Range("G27").NumberFormat = '#,##0'
$ws.Range("G27").Value = 1
$ws.Range("H27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H27").Value = -100
